$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Adalogger M0"
$ws.Range("A1").Select()
